$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($cellRef, $text)
    # Force literal text storage (avoid Excel auto-converting numeric-looking
    # strings like "59.90" or "1.00" into real numbers, which would drop
    # meaningful trailing zeros / alter the displayed value).
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText "D2" "41.935.07"
$ws.Range("E2").Value = "  +1.28%  "

# Row 3 - Ethereum
Set-PriceText "D3" "2.217.97"
$ws.Range("E3").Value = "  +1.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
Set-PriceText "D5" "251.33"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.43%  "

# Row 7 - Solana
Set-PriceText "D7" "67.95"
$ws.Range("E7").Value = "  -1.27%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - Cardano
Set-PriceText "D9" "0.632"
$ws.Range("E9").Value = "  +7.74%  "

# Row 10 - Avalanche
Set-PriceText "D10" "39.26"
$ws.Range("E10").Value = "  +3.52%  "

# Row 11 - OKB
Set-PriceText "D11" "59.90"
$ws.Range("E11").Value = "  +2.74%  "

# Row 12 - Dogecoin
Set-PriceText "D12" "0.0938"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13 - Polkadot
Set-PriceText "D13" "7.09"
$ws.Range("E13").Value = "  -1.35%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.92%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-PriceText "D15" "2.550.74"
$ws.Range("E15").Value = "  +1.01%  "

# Row 16 - Chainlink
$ws.Range("E16").Value = "  -0.53%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -0.89%  "

# Row 18 - WrappedEther
Set-PriceText "D18" "2.215.33"
$ws.Range("E18").Value = "  +0.60%  "

# Row 19 - WrappedBTC
Set-PriceText "D19" "41.859.52"
$ws.Range("E19").Value = "  +1.21%  "

# Row 20 - ShibaInu
Set-PriceText "D20" "0.0₃0961"
$ws.Range("E20").Value = "  +0.89%  "

# Row 21 - Uniswap
Set-PriceText "D21" "6.21"
$ws.Range("E21").Value = "  -0.68%  "

# Row 22 - Litecoin
Set-PriceText "D22" "72.62"
$ws.Range("E22").Value = "  +1.07%  "

# Row 23 - BitcoinCash (price only)
Set-PriceText "D23" "231.91"

# Row 24 - ImmutableX
Set-PriceText "D24" "2.06"
$ws.Range("E24").Value = "  -2.01%  "

# Row 25 - WEMIXToken
Set-PriceText "D25" "3.90"
$ws.Range("E25").Value = "  +0.53%  "

# Row 26 / 27 - Dai and Cosmos swap positions, with updated data
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-PriceText "D26" "11.40"
$ws.Range("E26").Value = "  -5.93%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-PriceText "D27" "1.00"
$ws.Range("E27").Value = "  +0.16%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -3.94%  "

# Row 29 - LEO
Set-PriceText "D29" "3.71"
$ws.Range("E29").Value = "  -1.19%  "

# Row 30 - Toncoin
Set-PriceText "D30" "2.25"
$ws.Range("E30").Value = "  +2.52%  "

# Row 31 - Monero
Set-PriceText "D31" "166.96"
$ws.Range("E31").Value = "  -1.81%  "

# Row 32 - EthereumClassic
Set-PriceText "D32" "20.44"
$ws.Range("E32").Value = "  -1.12%  "

# Row 33 - Hedera
Set-PriceText "D33" "0.0802"
$ws.Range("E33").Value = "  +9.78%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +6.45%  "

# Row 35 - Kaspa
Set-PriceText "D35" "0.121"
$ws.Range("E35").Value = "  -0.97%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  -0.18%  "

# Row 37 - Filecoin
Set-PriceText "D37" "4.62"
$ws.Range("E37").Value = "  -0.45%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +2.72%  "

# Row 39 - InjectiveProtocol
Set-PriceText "D39" "25.57"
$ws.Range("E39").Value = "  -3.13%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +2.44%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  +0.42%  "

# Row 42 - Celestia
Set-PriceText "D42" "12.20"
$ws.Range("E42").Value = "  +1.09%  "

# Row 43 - THORChain
Set-PriceText "D43" "5.66"
$ws.Range("E43").Value = "  -2.19%  "

# Row 44 - FTXToken
Set-PriceText "D44" "5.09"
$ws.Range("E44").Value = "  +1.37%  "

# Row 45 - MultiversX
Set-PriceText "D45" "62.08"
$ws.Range("E45").Value = "  -3.01%  "

# Row 46 - Algorand
Set-PriceText "D46" "0.199"
$ws.Range("E46").Value = "  -2.62%  "

# Row 47 - FraxShare
Set-PriceText "D47" "8.59"
$ws.Range("E47").Value = "  -0.69%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  -1.15%  "

# Row 49 - BinanceUSD
$ws.Range("E49").Value = "  -0.64%  "

# Row 50 - ARBITRUM
Set-PriceText "D50" "1.16"
$ws.Range("E50").Value = "  +1.03%  "

# Row 51 - SynthetixNetwork
Set-PriceText "D51" "4.36"
$ws.Range("E51").Value = "  +2.82%  "
